$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as text in the
# source workbook (t="inlineStr"). Force the number format to Text before
# assigning so Excel doesn't auto-convert these into numeric cells.
$priceCells = @(
    "D2","D3","D4","D6","D7","D9","D11","D12","D13","D14","D15","D16",
    "D17","D18","D19","D21","D22","D23","D24","D25","D40","D41","D42",
    "D43","D44","D45","D47","D48","D49"
)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Price-only updates
$ws.Range("D2").Value = "238.63"
$ws.Range("D3").Value = "21.89"
$ws.Range("D4").Value = "5.383"
$ws.Range("D6").Value = "6.475"
$ws.Range("D7").Value = "3.341"
$ws.Range("D9").Value = "1.026"
$ws.Range("D11").Value = "0.07339"
$ws.Range("D12").Value = "0.03156"
$ws.Range("D13").Value = "0.02978"
$ws.Range("D14").Value = "0.09234"
$ws.Range("D15").Value = "0.001660"
$ws.Range("D16").Value = "3.255"
$ws.Range("D17").Value = "0.04770"
$ws.Range("D18").Value = "0.0005714"
$ws.Range("D19").Value = "0.006256"
$ws.Range("D21").Value = "0.001052"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D23").Value = "0.0004003"
$ws.Range("D24").Value = "3.914"
$ws.Range("D25").Value = "2.201"
$ws.Range("D40").Value = "0.04078"
$ws.Range("D41").Value = "0.006961"

# Row 42 and 43 swap: CEJI and BKEXToken swap places, with new prices
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1039"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003002"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "0.009158"
$ws.Range("D45").Value = "0.00005439"
$ws.Range("D47").Value = "0.6757"
$ws.Range("D48").Value = "0.03760"
$ws.Range("D49").Value = "0.00002101"
